# Update cryptocurrency price/volume figures in the worksheet to reflect the
# latest scrape (GitHub Actions scheduled refresh).
#
# Column D ("Price") and column E ("Volume(1h)") hold values that are stored
# as *text* in the workbook (prices such as "28.057.41" use dots as
# thousands separators rather than being real numbers, and the volume
# column keeps its literal leading/trailing spaces). Several of the new
# price strings would otherwise be auto-recognised by Excel as numbers, so
# those cells are pre-formatted as Text before the value is written to make
# sure they keep behaving like the original inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '28.057.41'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -3.76%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.917.13'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -2.93%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -1.02%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '330.49'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +0.30%  '; ForceText = $false },
    @{ Cell = 'E6'; Value = '  -0.87%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.4697'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -5.38%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.4032'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -4.34%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '53.18'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -1.76%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '0.08395'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  -9.94%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '1.044'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -5.14%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '22.11'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  -3.30%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '1.908.89'; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -3.37%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '7.462'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -5.60%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '6.065'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -6.18%  '; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -1.05%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '89.89'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -2.35%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '0.00001060'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -4.77%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.06577'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -2.28%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '18.06'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -6.03%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '1.002'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -0.93%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '5.718'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -4.12%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '28.041.50'; ForceText = $false },
    @{ Cell = 'E23'; Value = '  -3.86%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '11.34'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -5.39%  '; ForceText = $false },
    @{ Cell = 'E25'; Value = '  +1.04%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '2.119.09'; ForceText = $false },
    @{ Cell = 'E26'; Value = '  -4.18%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '153.85'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -1.56%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '19.99'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -3.91%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '2.133'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -6.04%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '5.706'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -9.05%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '123.30'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -3.21%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '0.9716'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -7.20%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '0.09583'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -2.75%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '1.444'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -4.16%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  -2.58%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '5.529'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -5.03%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '8.917'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -1.50%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.02305'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -4.98%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.06136'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -4.61%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '1.220'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -8.26%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.6134'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -5.40%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '11.01'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -4.42%  '; ForceText = $false },
    @{ Cell = 'E43'; Value = '  -0.90%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '1.306'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -3.78%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.5866'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -5.76%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '12.77'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -4.14%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '2.022'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -7.40%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '3.480'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -0.09%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.06827'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -2.10%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '109.70'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -3.14%  '; ForceText = $false }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    if ($update.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $update.Value
}
